$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)
$ws1 = $wb.Worksheets.Item(1)

# --- Step 1: Add new nutrition data rows (137-143) to NutritionalData sheet ---
# Row 137: garbanzo beans
$ws3.Range("A137").Value = "garbanzo beans 1/2 cup "
$ws3.Range("B137").Value = 120
$ws3.Range("C137").Value = 2
$ws3.Range("D137").Value = 0
$ws3.Range("E137").Value = 7
$ws3.Range("F137").Value = 20
$ws3.Range("G137").Value = 4
$ws3.Range("H137").Value = 420

# Row 138: pickled beets canned
$ws3.Range("A138").Value = "pickled beets canned 1/2 cup serving"
$ws3.Range("B138").Value = 40
$ws3.Range("C138").Value = 0
$ws3.Range("D138").Value = 0
$ws3.Range("E138").Value = 1
$ws3.Range("F138").Value = 8
$ws3.Range("G138").Value = 1
$ws3.Range("H138").Value = 140

# Row 139: Extra Virgin Olive oil and Canola oil mix
$ws3.Range("A139").Value = "Extra Virgin Olive oil and Canola oil mix to fry falafels in 1 tbs serving"
$ws3.Range("B139").Value = 120
$ws3.Range("C139").Value = 14
$ws3.Range("D139").Value = 1
$ws3.Range("E139").Value = 0
$ws3.Range("F139").Value = 0
$ws3.Range("G139").Value = 0
$ws3.Range("H139").Value = 0

# Row 140: lemon juice dressing
$ws3.Range("A140").Value = "lemon juice for the olive oil and lemon juice dressing 1 tbs"
$ws3.Range("B140").Value = 2
$ws3.Range("C140").Value = 0
$ws3.Range("D140").Value = 0
$ws3.Range("E140").Value = 0
$ws3.Range("F140").Value = 0
$ws3.Range("G140").Value = 0
$ws3.Range("H140").Value = 25

# Row 141: red pepper dressing
$ws3.Range("A141").Value = "The red pepper dressing 2 tbs serving"
$ws3.Range("B141").Value = 180
$ws3.Range("C141").Value = 16
$ws3.Range("D141").Value = 3
$ws3.Range("E141").Value = 7
$ws3.Range("F141").Value = 3
$ws3.Range("G141").Value = 1
$ws3.Range("H141").Value = 70

# Row 142: falafel dry mix
$ws3.Range("A142").Value = "falafel dry mix at Hummus Republic Chino, serving is 2 oz they get fried made of chickpeas, fava beans, and spices"
$ws3.Range("B142").Value = 200
$ws3.Range("C142").Value = 2.5
$ws3.Range("D142").Value = 0
$ws3.Range("E142").Value = 14
$ws3.Range("F142").Value = 31
$ws3.Range("G142").Value = 13
$ws3.Range("H142").Value = 590

# Row 143: falafel and hummus bowl totals (formulas)
$ws3.Range("A143").Value = "falfel and hummus bowl with dressing from above ingredients Hummus Republic 3-8-2021 6 pm"
$ws3.Range("B143").Formula = "=SUM(B142*5,B141,B140,B139,B138/2,B137*2)"
$ws3.Range("C143").Formula = "=SUM(C142*5,C141,C140,C139,C138/2,C137*2)"
$ws3.Range("D143").Formula = "=SUM(D142*5,D141,D140,D139,D138/2,D137*2)"
$ws3.Range("E143").Formula = "=SUM(E142*5,E141,E140,E139,E138/2,E137*2)"
$ws3.Range("F143").Formula = "=SUM(F142*5,F141,F140,F139,F138/2,F137*2)"
$ws3.Range("G143").Formula = "=SUM(G142*5,G141,G140,G139,G138/2,G137*2)"
$ws3.Range("H143").Formula = "=SUM(H142*5,H141,H140,H139,H138/2,H137*2)"

# --- Step 2: Update researchMeasures row 54 (3-8-2021 diary entry) ---
$avocadoText = @"
3/4 avocado
(241.5	21.75	3	3	12.75	13.5	10.5)
2 eggs
(140	10	3	12	0	0	140)
2 tbs olive oil
(120	14	2	0	0	0	0)
mozzarella cheese 1/3 cup
(106.7  6.7	4.7	8	1.3	0	253.3)
4 corn tortillas Guerrero
(200	2	0	4	42	4	40)
2 tbs sourcream
(60	5	3.5	1	2	0	15)

pea protein
(130	2	0	18	9	2	320)
banana
(105	0	0	1	27	3	1)
honey crisp apple
(62	0.1	0	0.3	14.9	2.5	0)
2 tbs cocoa
(20	1	0	2	6	2	0)
handful cashews 1/3 cup approximately
(164	13.5	2.5	4.7	8.4	0.9	4)
1 cup almond milk
(30	2.5	0	1	1	0	115)
bowl falafel with 3 scoops of hummus differnt flavors, beets, and cucumbers and dressing
(1562	46.5	4	91.5	202	74.5	3955)
grapefruit
(92	0	0	2	24	2	0)

=241.5+140+120+106.7+200+60+130+105+62+20+164+30+92+1562+92
=21.75+10+14+6.7+2+5+2+0+0.1+1+13.5+2.5+0+46.5+0
=3+3+2+4.7+0+3.5+0+0+0+0+2.5+0+0+4+0
=3+12+0+8+4+1+18+1+0.3+2+4.7+1+2+91.5+2
=12.75+0+0+1.3+42+2+9+27+14.9+6+8.4+1+24+202+24
=13.5+0+0+0+4+0+2+3+2.5+2+0.9+0+2+74.5+2
=10.5+140+0+253.3+40+15+320+1+0+0+4+115+0+3955+0




"@

$diaryText = @"
Woke up at 3 am. Went to the kitchen the roommate had it blocked off from pups to cook and we argued about him wanting me out of his way as he always has to argue about something. I wasn't even in his way. His ignorant male and white privilege. His mom doesn't even want him around because he is ignorant and thinks people are supposed to do what he says and not argue. Such a loser. But any ways. I planned on getting up at 3 am because of my study plans, and went to bed early after not getting any sleep since 11 pm the day before while completing homework due yesterday after work. Made my coffee while he argued and bitched, had a lg BM slightly dehydrated at first around 3:50 am. Looked at the week 3 chemistry powerpoint and made some notes on my notecards, then made a 2nd cup of coffee after feeding the babies, didn't drink it, but took my measurements at 5 am after reviewing the slides on nomenclature. Took a nap at just before 6 am and woke up just before 630 am and stayed in bed until the alarm for 630 am went off. Made breakfast of 3 eggs and 2 tbs sourcream blended in my Ninga bullet blender and pan scrambled in 2 tbs olive oil, and the rest of the Guerrero corn tortillas, had to chop off 1/8 of all of them for the hard crust at the edge on the same side of all, with about 1/2 cup mozzarella cheese and paprika. Shared about a third of that with my babies, Also had a whole lg avocado with the scrambled eggs and quesadillas. I probably ate 2 1/4 quesadillas, 1/3 cup mozzarella cheese, 2 eggs, and 3/4 avocado. Did the dishes, fed the neighborhood cat, that wasn't waiting for me to feed her/him. Meow meow doesn't eat human food. Then started reading genetics ch 4 and looking over the chemistry chapter 4 powerpoint downloaded earlier, and planning to workout or not before work at 3 pm. The computer froze up again on allowing me to use my cursor to click on anything I moved it to and I shut down the computer before the nap of 20 minutes earlier. It might be mechanical, but it was after logging onto the course website host , Canvas. But I was able to move the cursor but not able to use the click to select links or areas I move the cursor to. It could be the laptop cover not fitting well after taking it off a few times in the past for freezing up and not shutting off. Also, its uneven by the mousepad and cracked on the plastic, flimsy case housing of the device. I have enough instant coffee for my 3rd cup of coffee, but I need more instant coffee and to pick up my Amazon package at the hub locker, because it arrived yesterday, my eye liner and mascara. I got my eye pillow with lavendar fennel seeds and removable satin eye pillow covers delivered at my house yesterday, and wasn't expecting them until last week. The driver left it right by the sidewalk and the roommate found it. Fortunately, he was doing yard work yesterday and saw it. It could have been taken by anybody. Especially if anybody would walk by it regularly and see the box everyday in the same spot for a week. Cloudy outside and it makes me tired. It was cloudy all day yesterday too. I read the first 4 sections of chapter 4 in genetics and took another 20 minute nap, but had my work clothes in the wash before taking a nap and spun them to get out excess moisture too. When i woke up the clothes in dryer were still very damp and I couldn't put my laundry in the dryer. I made a smoothie with the frozen banana I had to peel frozen and it didn't peel easily, also a honey crisp apple that I only bought 1 of at Sprouts when I went because I dont' really eat apples, a serving of pea protein or two scoops, 2 tbs cocoa and a handful of raw unroasted cashews. It came out thick. I also put about a cup of almond unsweetened milk the Silk brand in it.  Had to add the banana after mixing it with apple first because too cramped. I had to eat it with a spoon because it was too thick. I have been waiting to eat a grapefruit, I grabbed it and put it on the desk to eat but every time I look at it I don't want to eat it. I think I will work out. The sun is shining a little bit. I need to check my amazon locker on the way to work. I had a potential lymphatic drainage client call that wants to email me the downloaded consent form because her browser isn't letting her fill it out and send it in. Haven't received it, she called before my nap. After swithching out the laundry and reentering my nutrition data and such, I had another BM before working out. The BM was a constipated one. Might be from the protein shakes that have been thick lately, not drinking enough water. That and the waist trimmer, I put it on with my workout clothes before working out and after drinking the protein smoothie. My skin is itching at my lower legs and arms and all over. Maybe need more water. Could be stress.Wearing tights. Also, there more solid than before the protein smoothies. Because the texture before starting the protein shakes was more soft like an actual poop emoji, these are almost turds. I put my face routine of the witch hazel, the skin pigmentation corrector, the revitalift, the face sunscreen, then my makeup, and had some extra time. Felt more movement in lower abs, like more or another BM on the way. But feel like it will be another constipated one. The last one I had was between my face routine after putting on the revitalift (it was small, the one before regular sized and the first one was large) to push hard and hold my breath with the waist trimmer on that it gave me acid reflux and I vomited a little stomach acid. It didn't burn my esophagus thankfully, more like spittal. Was able to do my workout after my 2nd BM and had a little bit of burping and slight acid reflux but no spittal. It took an hour and started at 1110 am and ended about 1210 pm, 62 degrees when I ended the workout. Still cloudy but the sun popped up some, made me worry because I didn't have sunscreen on and might have a hormonal skin condition that makes my face get brown spots like my mom had around my age. I feel like I might have a 4th BM. Not sure. But have to check the amazon locker and see if time to get some instant coffee and coincidentally toilet paper. Before leaving for work I did have a 4th BM, very small the size of what was probably in my rectum that I couldn't squeeze out from last BM. The protein shakes/smoothies make stuff just push right through you, I guess. The digestive tract works by entering the ascending colon on the right after going through the small intestine through a sphincter, then moving right to left through the transverse colon, then through the descending colon, and then to the rectum where it awaits to be vacated. Some people who are dehydrated, stressed, and not operating at a healthy level their autonomous nervous system, the sympathetic part don't have regular BMs and it is a sign your healthy depending on consistency, frequency, changes to diet, health, life, etc like medications. The parasympathetic nervous system is the fight or flight part of autonomous nervous system that shuts down digestion and hence regular BMs, sleep patterns etc. This is why people stress eat without realizing it, to operate the sympathetic nervous system and not let the stress affect their immune systems. At work I wasn't feeling good again. I still felt like I had to crap, and when I got there I did once before the 1 st massage, the first time taking a dump at work. Then again after the 2nd massage, twice, making the total BM for the day 7. Not sure why, but last Tuesday I felt better after loosening my belt and my waist trimmer. I loosened my belt and felt much better, then on my break ate a bowl from Hummus Republic, 5 falafels, feta cheese about 1/4 cup, some lemon oil sauce and red pepper sauce, and the hummus flavors were a red pepper hummus called harrins or something similar, a zesty hummus, and eggplant hummus, with pickles beats cut into tiny cubes and cucumbers, when I asked for the nutrition facts because I am watching my diet, they didn't have the nutrition facts on a sheet like jamba juice, but took photos of the indredient box and jars for the garbanzo beans for the hummus, the falafels about 5 small sized falafels that could fit into 1 cup fried in oil, the photo showed canola oil, and a photo of dry ingredients for the falafel has a serving as 2 oz at 200 calories per 2 oz serving. Google says there are 8 ounces in a cup. So I had four servings of the falafel. It was good, but spicy. Made my sinuses run. I also ate a grapefruit. Before my next two clients I felt like I might have indigestion again and a BM so I took off my waist trimmer after the next one and put back on but on the 31" setting instead of 30" setting and loosened my belt. I felt better. It must be the waist trimmer squeezing my insides and making me get indigestion and BMs. Weird that it only happens on certain days. This could be a part of the next research once concluded with this data set. To see or answer the question of what was causing me flu like indigestion that disappeared after 20 minutes randomly occuring by triggers like scent or new changes to schedule. I should add more features to the data like if pea protein added protein powder taken, how many days since last workout of weight lifting and separately for cardio kickboxing, and others like if cocoa powder used or a new fruit , etc. . I also got a couple private clients, return this Friday and one new one tomorrow morning but not the same one who called earlier in the day to send me photos of the consent form. I picked up my amazon makeup on the way to work and had to reschedule that appointment from 4 pm to 11 am after talking with client for her bday. That could have stressed me out more than I thought and triggered the indigestion at work. I got a text while at work that my other package from Amazon arrived and forgot to pick up the aquagel that I ordered for the RF machine and is waiting to be picked up after work. A couple who let me have the 5 samples of CBD rebooked for this Friday after work during my shift and I didn't get the news until my break at around 530 pm. All my clients are great people. I have an 11 am tomorrow and chemistry at 3 pm, and lab at 6 pm that the procedures need to be written in the lab and turned in before lab tomorrow. There is a lot of sodium in the falafels and garbanzo beans, but also a lot of protein. But wow! the sodium for that bowl of falafel and hummus is a lot. No wonder my ankles get swollen after hummus, so much sodium.
"@

$ws1.Range("AA54").Value = $avocadoText
$ws1.Range("Z54").Value = $diaryText

$ws1.Range("AB54").Formula = "=241.5+140+120+106.7+200+60+130+105+62+20+164+30+92+1562+92"
$ws1.Range("AC54").Formula = "=21.75+10+14+6.7+2+5+2+0+0.1+1+13.5+2.5+0+46.5+0"
$ws1.Range("AD54").Formula = "=3+3+2+4.7+0+3.5+0+0+0+0+2.5+0+0+4+0"
$ws1.Range("AE54").Formula = "=3+12+0+8+4+1+18+1+0.3+2+4.7+1+2+91.5+2"
$ws1.Range("AF54").Formula = "=12.75+0+0+1.3+42+2+9+27+14.9+6+8.4+1+24+202+24"
$ws1.Range("AG54").Formula = "=13.5+0+0+0+4+0+2+3+2.5+2+0.9+0+2+74.5+2"
$ws1.Range("AH54").Formula = "=10.5+140+0+253.3+40+15+320+1+0+0+4+115+0+3955+0"

$ws1.Range("AP54").Value = 7

# --- Step 3: Update view/selection state ---
$ws3.Activate()
$ws3.Range("H143").Select()
$ws1.Activate()
$ws1.Range("Z54").Select()
